$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Remove the existing hyperlink on C2 (the mailto:Tetherfi@930 link) before changing values
if ($ws.Range("C2").Hyperlinks.Count -gt 0) {
    $ws.Range("C2").Hyperlinks.Delete()
}

# Update row 2 values (A2, B2, C2) with new credentials / url
$ws.Range("A2").Value = "http://10.133.146.17:56080/SG/UOB_OCM"
$ws.Range("B2").Value = "meghna"
$ws.Range("C2").Value = "P@ssw0rd@123"

# Update the active selection to B3 (as in the edited workbook)
$ws.Range("B3").Select()
